# Update gh-pages output (苏州-漫展信息.xlsx) to the newly generated data.
# The same table lives on two sheets: "展览" (sheet 1) and "全部类型" (sheet 4).
# Both must receive identical updates.

$wb = $excel.ActiveWorkbook

foreach ($sheetIndex in 1, 4) {
    $ws = $wb.Worksheets.Item($sheetIndex)

    # --- simple "想去人数" (interest count) refreshes ---
    $ws.Range("F2").Value = 626
    $ws.Range("F4").Value = 81
    $ws.Range("F5").Value = 12968
    $ws.Range("F6").Value = 70
    $ws.Range("F10").Value = 1169
    $ws.Range("F12").Value = 13730
    $ws.Range("F13").Value = 14250
    $ws.Range("F22").Value = 1084
    $ws.Range("F28").Value = 290

    # --- rows 25 & 26: the two OCG events swap order/content, and the ---
    # --- "general" session's interest count is refreshed to 5327.    ---
    $ws.Range("C25").Value = "苏州·OCG国潮动漫游戏嘉年华"
    $ws.Range("E25").Value = "2024.05.04 09:00-05.05 17:00"
    $ws.Range("F25").Value = 5327
    $ws.Range("G25").Value = 65
    $ws.Range("H25").Value = "https://show.bilibili.com/platform/detail.html?id=82779"
    $ws.Range("I25").Value = "//i1.hdslb.com/bfs/openplatform/202403/hcgdIzw61710298907237.jpeg"

    $ws.Range("C26").Value = "苏州·OCG国潮动漫游戏嘉年华阿杰内场"
    $ws.Range("E26").Value = "2024.05.04 09:00-05.04 17:00"
    $ws.Range("F26").Value = 937
    $ws.Range("G26").Value = "已售罄"
    $ws.Range("H26").Value = "https://show.bilibili.com/platform/detail.html?id=82940"
    $ws.Range("I26").Value = "//i2.hdslb.com/bfs/openplatform/202403/lLKmv48C1710511298160.jpeg"
}
